# Rename the three header/footer logo pictures (Word COM automation only
# exposes a rename operation on InlineShape via the .Name property - it
# writes through to the picture's <wp:docPr name="..."/> attribute).
#
#   footer (first page)  -> Pearson logo: image1.png -> image2.png
#   footer (default)     -> Pearson logo: image1.png -> image2.png
#   header (first page)  -> BTEC logo   : image2.jpg -> image1.jpg

$d = $word.ActiveDocument
$section = $d.Sections.Item(1)

# wdHeaderFooterPrimary = 1, wdHeaderFooterFirstPage = 2
$wdHeaderFooterPrimary = 1
$wdHeaderFooterFirstPage = 2

# --- Footers: Pearson logo, image1.png -> image2.png -------------------
$firstPageFooter = $section.Footers.Item($wdHeaderFooterFirstPage)
if ($firstPageFooter.Exists) {
    for ($i = 1; $i -le $firstPageFooter.Range.InlineShapes.Count; $i++) {
        $shape = $firstPageFooter.Range.InlineShapes.Item($i)
        $shape.Name = "image2.png"
    }
}

$defaultFooter = $section.Footers.Item($wdHeaderFooterPrimary)
if ($defaultFooter.Exists) {
    for ($i = 1; $i -le $defaultFooter.Range.InlineShapes.Count; $i++) {
        $shape = $defaultFooter.Range.InlineShapes.Item($i)
        $shape.Name = "image2.png"
    }
}

# --- Header: BTEC logo, image2.jpg -> image1.jpg ------------------------
$firstPageHeader = $section.Headers.Item($wdHeaderFooterFirstPage)
if ($firstPageHeader.Exists) {
    for ($i = 1; $i -le $firstPageHeader.Range.InlineShapes.Count; $i++) {
        $shape = $firstPageHeader.Range.InlineShapes.Item($i)
        $shape.Name = "image1.jpg"
    }
}
